$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (School ID 222444) was "West National High School" / "private".
# Update the school name to the corrected full name, and the type to "public".
$ws.Range("B9").Value = "Canumay West National High School"
$ws.Range("C9").Value = "public"
